# Auto-generated edit script: apply the crypto price/volume update from the commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.019.94'
$ws.Range("E2").Value = '  -1.61%  '

$ws.Range("D3").Value = '1.857.64'
$ws.Range("E3").Value = '  -2.82%  '

$ws.Range("E4").Value = '  +0.25%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.02'
$ws.Range("E5").Value = '  -2.53%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.14%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4651'
$ws.Range("E7").Value = '  -2.76%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2806'
$ws.Range("E8").Value = '  -1.38%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06520'
$ws.Range("E9").Value = '  -2.77%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.86'
$ws.Range("E10").Value = '  +4.59%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07739'
$ws.Range("E11").Value = '  +0.32%  '

$ws.Range("B12").Value = 'Litecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '95.96'
$ws.Range("E12").Value = '  -6.63%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.875.91'
$ws.Range("E13").Value = '  -1.92%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.036'
$ws.Range("E14").Value = '  -3.30%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6654'
$ws.Range("E15").Value = '  -0.88%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '280.17'
$ws.Range("E16").Value = '  +3.15%  '

$ws.Range("D17").Value = '30.079.92'
$ws.Range("E17").Value = '  -1.41%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.002'

$ws.Range("D19").Value = '2.148.22'
$ws.Range("E19").Value = '  -0.47%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.44'
$ws.Range("E20").Value = '  -1.87%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.301'
$ws.Range("E21").Value = '  -2.25%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.13%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.000007182'
$ws.Range("E23").Value = '  -3.89%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.115'
$ws.Range("E24").Value = '  -3.05%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '166.16'
$ws.Range("E25").Value = '  -0.39%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.239'
$ws.Range("E26").Value = '  -1.80%  '

$ws.Range("E27").Value = '  -2.78%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.927'
$ws.Range("E28").Value = '  -6.77%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.363'
$ws.Range("E29").Value = '  -1.40%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09618'
$ws.Range("E30").Value = '  -4.37%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.352'
$ws.Range("E31").Value = '  -5.51%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.461'
$ws.Range("E32").Value = '  -3.50%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.079'
$ws.Range("E33").Value = '  -3.90%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04638'
$ws.Range("E34").Value = '  -1.95%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6948'
$ws.Range("E35").Value = '  -4.60%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.079'
$ws.Range("E36").Value = '  -2.98%  '

$ws.Range("E37").Value = '  +0.25%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.700'
$ws.Range("E38").Value = '  -0.71%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01846'
$ws.Range("E39").Value = '  -3.82%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.329'
$ws.Range("E40").Value = '  +1.04%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.502'
$ws.Range("E41").Value = '  -4.12%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '71.00'
$ws.Range("E42").Value = '  -5.09%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8544'
$ws.Range("E43").Value = '  -0.56%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.924'
$ws.Range("E44").Value = '  -2.06%  '

$ws.Range("E45").Value = '  +0.12%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.90'
$ws.Range("E46").Value = '  -2.05%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4130'
$ws.Range("E47").Value = '  -3.15%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '968.57'
$ws.Range("E48").Value = '  +5.57%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.142'
$ws.Range("E49").Value = '  -3.84%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.048'
$ws.Range("E50").Value = '  +2.35%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '33.65'
$ws.Range("E51").Value = '  -3.37%  '
